# In database mode statistics results are also stored in the db.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 2: From-date changes, Interval changes from 5m to 30m
$ws.Range("D2").Value = 44501
$ws.Range("F2").Value = "30m"

# Rows 3 and 4 become new test-case rows; copy row 2's formatting down first
$ws.Range("A2:J2").Copy()
$ws.Range("A3:J4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 3: fill in a new full test-case row
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Bybit"
$ws.Range("C3").Value = "BTCUSDT"
$ws.Range("D3").Value = 44501
$ws.Range("E3").Value = 44561
$ws.Range("F3").Value = "30m"
$ws.Range("G3").Value = 10000
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = "ScalpEmaRsiAdx_X"

# Row 4: fill in a new full test-case row
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Bybit"
$ws.Range("C4").Value = "BTCUSDT"
$ws.Range("D4").Value = 44501
$ws.Range("E4").Value = 44561
$ws.Range("F4").Value = "30m"
$ws.Range("G4").Value = 10000
$ws.Range("H4").Value = 0.7
$ws.Range("I4").Value = 0.7
$ws.Range("J4").Value = "ScalpEmaRsiAdx_X"

# Update the active cell selection
$ws.Range("H5").Select()
